$wb = $excel.ActiveWorkbook

# --- Work on the "Repayment schedule" sheet ---
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (pushes old N/O/P -> O/P/Q)
$mWidth = $ws.Columns("M:M").ColumnWidth
$ws.Columns("N:N").Insert()

# Give the freshly inserted column N the same width as column M (~10.71)
$ws.Columns("N:N").ColumnWidth = $mWidth

# Update the selection on this sheet to K15
$selected = $ws.Range("K15").Select()

# Make "Repayment schedule" the active/selected sheet (activeTab = 2, 0-indexed)
$ws.Activate()

# --- "Transactions" sheet loses the tab-selected flag (handled automatically
#     since only one sheet can be the active tab at a time; activating the
#     Repayment schedule sheet above clears it from Transactions) ---

$saved = $wb.Save()
